$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "AAAKN3DD&E"
$ws.Range("A3").Value = "AAAKN3DD&E"
$ws.Range("A4").Value = "AAAKN3DD&E"
$ws.Range("A5").Value = "AAAKN3DD&E"

$ws.Range("B11").Select()
